$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64, shifting rows 64:89 down to 65:90
$ws.Rows.Item(64).Insert()

# Populate the new row 64 with the new data
$ws.Cells.Item(64, 1).Value = 4
$ws.Cells.Item(64, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(64, 3).Value = "Los Lagos"
$ws.Cells.Item(64, 4).Value = 44572
$ws.Cells.Item(64, 4).NumberFormat = $ws.Cells.Item(65, 4).NumberFormat
$ws.Cells.Item(64, 5).Value = 10
$ws.Cells.Item(64, 6).Value = 100112052
$ws.Cells.Item(64, 7).Value = "Albahaca"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 120
$ws.Cells.Item(64, 11).Value = 7000
$ws.Cells.Item(64, 12).Value = 7000
$ws.Cells.Item(64, 13).Value = 7000
$ws.Cells.Item(64, 14).Value = "$/docena de matas"
$ws.Cells.Item(64, 15).Value = "Región Metropolitana"
$ws.Cells.Item(64, 16).Value = 1167
$ws.Cells.Item(64, 17).Value = 6
$ws.Cells.Item(64, 18).Value = "Hortaliza"
